# PSXSymbols.xlsx — "removed misc scripts folder"
# Adds the KSE100/KSE30 index rows to MYLIST, drops the now-unwanted ALLSHR
# row from KMI30 / KMIALL (shifting KSE30 up), tags the index rows on CUSTUM
# with a monthly ("M") marker in column B, and leaves KMI30 as the active tab.

$wb = $excel.ActiveWorkbook

# --- MYLIST sheet: add two new rows (KSE100 / KSE30) ---
$wsMyList = $wb.Worksheets.Item("MYLIST")
$wsMyList.Range("A33").Value = "KSE100"
$wsMyList.Range("A34").Value = "KSE30"
$null = $wsMyList.Activate()
$null = $wsMyList.Range("A33:A34").Select()

# --- KMI30 sheet: remove the ALLSHR row (row 32), shifting KSE30 up ---
$wsKmi30 = $wb.Worksheets.Item("KMI30")
$null = $wsKmi30.Rows.Item(32).Delete()

# --- KMIALL sheet: remove the ALLSHR row (row 213), shifting KSE30 up ---
$wsKmiAll = $wb.Worksheets.Item("KMIALL")
$null = $wsKmiAll.Rows.Item(213).Delete()
$null = $wsKmiAll.Activate()
$null = $wsKmiAll.Range("A213:XFD213").Select()

# --- CUSTUM sheet: add "M" markers next to the three index rows ---
$wsCustum = $wb.Worksheets.Item("CUSTUM")
$wsCustum.Range("B23").Value = "M"
$wsCustum.Range("B24").Value = "M"
$wsCustum.Range("B25").Value = "M"
$null = $wsCustum.Activate()
$null = $wsCustum.Range("A2:A19").Select()

# --- KMI30 sheet becomes the final active tab, selection on E31 ---
$null = $wsKmi30.Activate()
$null = $wsKmi30.Range("E31").Select()

Write-Output "Edit applied."
